$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells F2:F4 are formatted as Text ("@"), which would normally coerce any
# assigned number into a text/string cell. To write genuine numeric phone
# values while preserving the original Text cell format/style, temporarily
# switch the format to General, write the number, then restore the Text
# format.
$cells = @("F2", "F3", "F4")
$values = @(15911111111, 15911111112, 15911111113)

for ($i = 0; $i -lt $cells.Length; $i++) {
    $cell = $ws.Range($cells[$i])
    $originalFormat = $cell.NumberFormat
    $cell.NumberFormat = "general"
    $cell.Value = $values[$i]
    $cell.NumberFormat = $originalFormat
}

$ws.Range("F2").Select()
